$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 26,6

$arr[0,0] = 0.6753301551942219
$arr[0,1] = 1.667794583268128
$arr[0,2] = 0.8054896365839992
$arr[0,3] = 0.496779210170732
$arr[0,4] = 1
$arr[0,5] = 3.645393585217082

$arr[1,0] = 0.6753301551942219
$arr[1,1] = 1.667794583268128
$arr[1,2] = 0.1575252929769615
$arr[1,3] = 0.496779210170732
$arr[1,4] = 0
$arr[1,5] = 2.997429241610044

$arr[2,0] = 1.459612070389937
$arr[2,1] = 1.667794583268128
$arr[2,2] = 0.1575252929769615
$arr[2,3] = 0.496779210170732
$arr[2,4] = 1
$arr[2,5] = 3.781711156805759

$arr[3,0] = 0.3048080303191223
$arr[3,1] = 0.002777888934908601
$arr[3,2] = 0.1575252929769615
$arr[3,3] = 0.496779210170732
$arr[3,4] = 0
$arr[3,5] = 0.9618904224017244

$arr[4,0] = 3.230985683306322
$arr[4,1] = 1.667794583268128
$arr[4,2] = 0.8054896365839992
$arr[4,3] = 0.496779210170732
$arr[4,4] = 0
$arr[4,5] = 6.201049113329182

$arr[5,0] = 3.230985683306322
$arr[5,1] = 1.667794583268128
$arr[5,2] = 0.1575252929769615
$arr[5,3] = 0.496779210170732
$arr[5,4] = 0
$arr[5,5] = 5.553084769722144

$arr[6,0] = 0.3048080303191223
$arr[6,1] = 0.04240448674262143
$arr[6,2] = 0.1575252929769615
$arr[6,3] = 0.496779210170732
$arr[6,4] = 1
$arr[6,5] = 1.001517020209437

$arr[7,0] = 0.6753301551942219
$arr[7,1] = 1.667794583268128
$arr[7,2] = 0.8054896365839992
$arr[7,3] = 0.496779210170732
$arr[7,4] = 1
$arr[7,5] = 3.645393585217082

$arr[8,0] = 3.230985683306322
$arr[8,1] = 1.667794583268128
$arr[8,2] = 0.8054896365839992
$arr[8,3] = 8.660232485948974
$arr[8,4] = 1
$arr[8,5] = 14.36450238910742

$arr[9,0] = 0.127881588408715
$arr[9,1] = 0.04240448674262143
$arr[9,2] = 0.8054896365839992
$arr[9,3] = 0.496779210170732
$arr[9,4] = 0
$arr[9,5] = 1.472554921906068

$arr[10,0] = 3.230985683306322
$arr[10,1] = 10.29869402782916
$arr[10,2] = 0.8054896365839992
$arr[10,3] = 8.660232485948974
$arr[10,4] = 0
$arr[10,5] = 22.99540183366846

$arr[11,0] = 0.6753301551942219
$arr[11,1] = 1.667794583268128
$arr[11,2] = 0.8054896365839992
$arr[11,3] = 0.496779210170732
$arr[11,4] = 0
$arr[11,5] = 3.645393585217082

$arr[12,0] = 3.230985683306322
$arr[12,1] = 1.667794583268128
$arr[12,2] = 3.900430680208489
$arr[12,3] = 8.660232485948974
$arr[12,4] = 1
$arr[12,5] = 17.45944343273191

$arr[13,0] = 0.01514828764759746
$arr[13,1] = 0.04240448674262143
$arr[13,2] = 0.1575252929769615
$arr[13,3] = 0.496779210170732
$arr[13,4] = 0
$arr[13,5] = 0.7118572775379124

$arr[14,0] = 3.230985683306322
$arr[14,1] = 1.667794583268128
$arr[14,2] = 3.900430680208489
$arr[14,3] = 8.660232485948974
$arr[14,4] = 1
$arr[14,5] = 17.45944343273191

$arr[15,0] = 0.3048080303191223
$arr[15,1] = 0.3127903958511391
$arr[15,2] = 3.900430680208489
$arr[15,3] = 0.496779210170732
$arr[15,4] = 1
$arr[15,5] = 5.014808316549482

$arr[16,0] = 1.459612070389937
$arr[16,1] = 1.667794583268128
$arr[16,2] = 0.8054896365839992
$arr[16,3] = 645.32727682996
$arr[16,4] = 1
$arr[16,5] = 649.2601731202021

$arr[17,0] = 0.04763786555579896
$arr[17,1] = 0.3127903958511391
$arr[17,2] = 0.8054896365839992
$arr[17,3] = 0.496779210170732
$arr[17,4] = 0
$arr[17,5] = 1.662697108161669

$arr[18,0] = 0.6753301551942219
$arr[18,1] = 0.3127903958511391
$arr[18,2] = 0.8054896365839992
$arr[18,3] = 0.496779210170732
$arr[18,4] = 0
$arr[18,5] = 2.290389397800092

$arr[19,0] = 3.230985683306322
$arr[19,1] = 1.667794583268128
$arr[19,2] = 0.1575252929769615
$arr[19,3] = 0.496779210170732
$arr[19,4] = 1
$arr[19,5] = 5.553084769722144

$arr[20,0] = 0.6753301551942219
$arr[20,1] = 1.667794583268128
$arr[20,2] = 0.1575252929769615
$arr[20,3] = 0.496779210170732
$arr[20,4] = 1
$arr[20,5] = 2.997429241610044

$arr[21,0] = 3.230985683306322
$arr[21,1] = 1.667794583268128
$arr[21,2] = 0.1575252929769615
$arr[21,3] = 0.496779210170732
$arr[21,4] = 1
$arr[21,5] = 5.553084769722144

$arr[22,0] = 3.230985683306322
$arr[22,1] = 10.29869402782916
$arr[22,2] = 3.900430680208489
$arr[22,3] = 8.660232485948974
$arr[22,4] = 1
$arr[22,5] = 26.09034287729295

$arr[23,0] = 0.127881588408715
$arr[23,1] = 0.04240448674262143
$arr[23,2] = 3.900430680208489
$arr[23,3] = 0.496779210170732
$arr[23,4] = 1
$arr[23,5] = 4.567495965530558

$arr[24,0] = 0.6753301551942219
$arr[24,1] = 1.667794583268128
$arr[24,2] = 0.1575252929769615
$arr[24,3] = 0.496779210170732
$arr[24,4] = 0
$arr[24,5] = 2.997429241610044

$arr[25,0] = 3.230985683306322
$arr[25,1] = 0.3127903958511391
$arr[25,2] = 0.1575252929769615
$arr[25,3] = 0.496779210170732
$arr[25,4] = 0
$arr[25,5] = 4.198080582305154

$ws.Range("B2:G27").Value = $arr

Write-Host "done"
